# Update Name of Algo
# Applies updated RandomForest imputation result values to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -11.7504
$ws.Range("B3").Value = 6.017000000000001
$ws.Range("D3").Value = -6.930399999999995
$ws.Range("D12").Value = -7.341000000000005
$ws.Range("B14").Value = 5.799000000000003
$ws.Range("B16").Value = 6.037700000000002
$ws.Range("C18").Value = -12.1597
$ws.Range("B21").Value = 9.570299999999998
$ws.Range("B23").Value = 8.980200000000009
$ws.Range("C24").Value = -12.76299999999999
$ws.Range("D24").Value = -7.803799999999999
$ws.Range("B25").Value = 5.391600000000003
$ws.Range("C25").Value = -12.9335
$ws.Range("D25").Value = -8.69419999999999
$ws.Range("B26").Value = 6.145100000000002
$ws.Range("C27").Value = -13.1607
$ws.Range("B29").Value = 5.105500000000004
$ws.Range("C30").Value = -12.70189999999999
$ws.Range("C31").Value = -13.29589999999999
$ws.Range("C39").Value = -12.3468
$ws.Range("B40").Value = 9.043199999999995
$ws.Range("D41").Value = -7.967299999999995
$ws.Range("C42").Value = -12.6342
$ws.Range("C48").Value = -11.7131
$ws.Range("D50").Value = -8.309600000000001
$ws.Range("C51").Value = -11.5102
$ws.Range("C52").Value = -11.44689999999999
$ws.Range("B53").Value = 5.189699999999999
$ws.Range("D53").Value = -5.979800000000004
$ws.Range("C55").Value = -13.6077
$ws.Range("C56").Value = -12.6509
$ws.Range("D56").Value = -8.324400000000002
$ws.Range("B57").Value = 5.0928
$ws.Range("C57").Value = -13.08439999999999
$ws.Range("D57").Value = -8.723099999999995
$ws.Range("D58").Value = -8.185100000000004
$ws.Range("B59").Value = 4.768499999999996
$ws.Range("C60").Value = -13.76069999999999
$ws.Range("D61").Value = -7.940799999999997
$ws.Range("D63").Value = -7.931100000000004
$ws.Range("D64").Value = -7.9337
$ws.Range("B65").Value = 5.701700000000003
$ws.Range("B69").Value = 5.566299999999995
$ws.Range("D70").Value = -7.779799999999995
$ws.Range("D72").Value = -7.388199999999999
$ws.Range("C73").Value = -12.9928
$ws.Range("C74").Value = -12.3507
$ws.Range("B79").Value = 8.821800000000003
$ws.Range("B83").Value = 5.514599999999999
$ws.Range("D86").Value = -7.950499999999997
$ws.Range("C89").Value = -10.64950000000001
$ws.Range("D89").Value = -5.713000000000003
$ws.Range("C90").Value = -12.6786
$ws.Range("B91").Value = 4.953099999999997
$ws.Range("C92").Value = -11.3271
$ws.Range("B93").Value = 6.101800000000003
$ws.Range("D98").Value = -8.469900000000001
$ws.Range("B100").Value = 4.7918
$ws.Range("D100").Value = -8.378900000000005
$ws.Range("D102").Value = -7.938999999999993
